# Update calibrated cost values for IPPU prodinit rows.
# Columns J (10) through AS (45) hold the same value repeated across the
# time-series span for a given row; only the row-level constant changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    100 = 2495428.178
    101 = 1979420.265
    102 = 1945979.072
    103 = 2622894.833
    104 = 1919560.137
    105 = 1778429.437
    106 = 214057.6979
    107 = 7962037.544
    114 = 6045.46759
    115 = 7448734.738
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $ws.Range($ws.Cells.Item($row, 10), $ws.Cells.Item($row, 45)).Value = $value
}
